# Automatische test-sync: 2025-08-14 20:26:50
# Appends the new "Demo inplannen" log entry (row 5) to the Logs sheet,
# extends the conditional formatting ranges to cover the new row, and
# bumps the Dashboard summary count for the matching category.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 5 to the Logs sheet -------------------------------
$logs.Cells.Item(5, 1).Value = "Demo inplannen"
$logs.Cells.Item(5, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(5, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(5, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(5, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item(5, 6).Value = "2025-08-14 20:26:24"
$logs.Cells.Item(5, 7).Value = "Nee"
$logs.Cells.Item(5, 8).Value = "Ja"
$logs.Cells.Item(5, 9).Value = "Nee"
$logs.Cells.Item(5, 10).Value = "Nee"

# --- Extend the conditional formatting ranges from row 4 to row 5 -----
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# --- Update the Dashboard summary count --------------------------------
$dashboard.Range("B2").Value = 4
